# Updated data to reflect new requirement separation
# The "Prerequisites" table on the ART sheet gets three new columns inserted
# after "Prerequisites" (Corequisites, Concurrent, Recommended), pushing the
# old "Terms Typically Offered" column from D to G. Rows whose prerequisite
# text embedded a "Recommended: ..." clause have that clause split out into
# the new "Recommended" column (F) and removed from "Prerequisites" (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns at D (Corequisites, Concurrent, Recommended);
# this shifts the existing "Terms Typically Offered" column from D -> G,
# carrying all of its data with it, and grows dimension to A1:G79.
$ws.Range("D1:F1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Default all data rows (2-79) of the two brand-new columns to "NA"
$ws.Range("D2:E79").Value = "NA"

# Default the Recommended column to "NA" for every row; the handful of
# rows that actually had an embedded "Recommended: ..." clause are
# overwritten below with the real value.
$ws.Range("F2:F79").Value = "NA"

# --- Rows whose Prerequisites text said "one of the following:" -> "one of the" ---
$ws.Range("C31").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one of the ART 111, ART 112, ART 211, or ART 212."
$ws.Range("C32").Value = "One of the ART 112, ART 211, ART 212, or ART 213."
$ws.Range("C35").Value = "One of the ART 111, ART 112, ART 211, or ART 212."
$ws.Range("C52").Value = "Junior standing; and one of the ART 112, ART 212, ART 213, or WGS 201."

# --- Rows whose Prerequisites text had a trailing "Recommended: ..." clause ---
# split into column F, with the Prerequisites text trimmed and the Terms
# Typically Offered (now column G) value regenerated with a trailing space,
# matching the source data exactly.
$ws.Range("C37").Value = "ART 122 or ART 224."
$ws.Range("F37").Value = "ART 222."
$ws.Range("G37").Value = "F "

$ws.Range("C53").Value = "ART 182 or CSC 123."
$ws.Range("F53").Value = "ART 384."
$ws.Range("G53").Value = "W "

$ws.Range("C72").Value = "ART 122 or ART 182."
$ws.Range("F72").Value = "ART 384."
$ws.Range("G72").Value = "SP "

$ws.Range("C75").Value = "Junior standing."
$ws.Range("F75").Value = "ART 373 and ART 383."
$ws.Range("G75").Value = "SP "
